# PWMCalculations.xlsx update script
# - Swap/update the two "Field Weakening" / "Hall counter offset" shared strings
# - Bump the average Field-Weakening offset (F7) from 33 to 33.5 on the "Graph" sheet
# - Remove the now-obsolete helper calculations in I48/I49/J49
# - Make "Graph" the active/selected sheet (cell E7) instead of "Ref. Values"

$wb = $excel.ActiveWorkbook

$wsRef   = $wb.Worksheets.Item("Ref. Values")
$wsGraph = $wb.Worksheets.Item("Graph")

# --- Update the two related text notes (their roles are swapped) ---
$wsGraph.Range("A3").Value = "The offset added to the Hall counter is 23 for the states whit a falling edge of the Hall value and 44 to the state with a rising edge (avg = 33,5)"
$wsGraph.Range("D5").Value = "Field Weakening offset max"

# --- Update the average Field Weakening offset value, which drives the F/G columns (and the chart) ---
$wsGraph.Range("F7").Value = 33.5

# --- Remove the obsolete helper formulas that are no longer needed ---
$wsGraph.Range("I48").ClearContents()
$wsGraph.Range("I49").ClearContents()
$wsGraph.Range("J49").ClearContents()

# --- Switch the active tab/selection from "Ref. Values" to "Graph" ---
$wsRef.Range("C18:H18").Select()
$wsGraph.Activate()
$wsGraph.Range("E7").Select()

Write-Output "done"
